# The commit removes the slide that was just added (slide6 / sldId 2124817337,
# the "Stage_name" / "Short_quote" journey-stage slide) together with its
# notes page. Deleting the slide through the Slides collection removes the
# slide part, its notes slide part, and all of the now-dangling relationship
# entries (presentation.xml's sldIdLst entry, slideX.xml.rels, etc.) the same
# way PowerPoint itself does when a slide is deleted from the deck.

$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 2124817337) {
        $target = $slide
        break
    }
}

if ($target -eq $null) {
    # Fallback: the authored deck always had this slide last (slide 6).
    $target = $p.Slides.Item($p.Slides.Count)
}

$target.Delete()
